# Auto-generated edit script applying the scheduled-runner price updates
# to the Leve profit-tracking columns (H-N) across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 238.23077
$ws.Range("J2").Value = 226.11111
$ws.Range("L2").Value = 226.11111
$ws.Range("N2").Value = -452.11111
$ws.Range("H17").Value = 668.1900000000001
$ws.Range("J17").Value = 668.1900000000001
$ws.Range("L17").Value = 2004.57
$ws.Range("N17").Value = -2340.57
$ws.Range("H58").Value = 9239.333000000001
$ws.Range("J58").Value = 11340.833
$ws.Range("L58").Value = 34022.499
$ws.Range("N58").Value = -34322.499
$ws.Range("H107").Value = 3122.6667
$ws.Range("I107").Value = 1782.1428
$ws.Range("K107").Value = 1782.1428
$ws.Range("M107").Value = 137.8571999999999
$ws.Range("H115").Value = 358.6
$ws.Range("I115").Value = 358.6
$ws.Range("K115").Value = 1075.8
$ws.Range("M115").Value = 491.1999999999998
$ws.Range("H137").Value = 1572.85
$ws.Range("I137").Value = 936.8461
$ws.Range("J137").Value = 2754
$ws.Range("K137").Value = 2810.5383
$ws.Range("L137").Value = 8262
$ws.Range("M137").Value = -260.5383000000002
$ws.Range("N137").Value = -13362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9996.030000000001
$ws.Range("I61").Value = 8292.044
$ws.Range("K61").Value = 8292.044
$ws.Range("M61").Value = -8080.044
$ws.Range("H88").Value = 5157.1724
$ws.Range("I88").Value = 1400.8889
$ws.Range("J88").Value = 6847.5
$ws.Range("K88").Value = 1400.8889
$ws.Range("L88").Value = 6847.5
$ws.Range("M88").Value = -994.8888999999999
$ws.Range("N88").Value = -7659.5
$ws.Range("H91").Value = 5157.1724
$ws.Range("I91").Value = 1400.8889
$ws.Range("J91").Value = 6847.5
$ws.Range("K91").Value = 1400.8889
$ws.Range("L91").Value = 6847.5
$ws.Range("M91").Value = 3.111100000000079
$ws.Range("N91").Value = -9655.5
$ws.Range("H102").Value = 2452.5715
$ws.Range("I102").Value = 2378
$ws.Range("K102").Value = 2378
$ws.Range("M102").Value = -756
$ws.Range("H136").Value = 9996.030000000001
$ws.Range("I136").Value = 8292.044
$ws.Range("K136").Value = 24876.132
$ws.Range("M136").Value = -22326.132

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 99780
$ws.Range("J59").Value = 99780
$ws.Range("L59").Value = 99780
$ws.Range("N59").Value = -101474
$ws.Range("H105").Value = 2616.0557
$ws.Range("I105").Value = 2345.3076
$ws.Range("K105").Value = 2345.3076
$ws.Range("M105").Value = -598.3076000000001
$ws.Range("H107").Value = 2192
$ws.Range("I107").Value = 2182.25
$ws.Range("J107").Value = 2270
$ws.Range("K107").Value = 2182.25
$ws.Range("L107").Value = 2270
$ws.Range("M107").Value = -262.25
$ws.Range("N107").Value = -6110

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 564.72974
$ws.Range("I22").Value = 232.89473
$ws.Range("J22").Value = 915
$ws.Range("K22").Value = 232.89473
$ws.Range("L22").Value = 915
$ws.Range("M22").Value = 117.10527
$ws.Range("N22").Value = -1615
$ws.Range("H58").Value = 4532.091
$ws.Range("I58").Value = 3726.6365
$ws.Range("K58").Value = 3726.6365
$ws.Range("M58").Value = -3523.6365
$ws.Range("H62").Value = 13894.333
$ws.Range("I62").Value = 5398.4
$ws.Range("K62").Value = 5398.4
$ws.Range("M62").Value = -4774.4
$ws.Range("H65").Value = 13894.333
$ws.Range("I65").Value = 5398.4
$ws.Range("K65").Value = 26992
$ws.Range("M65").Value = -23872
$ws.Range("H132").Value = 2509.9167
$ws.Range("I132").Value = 2556.318
$ws.Range("K132").Value = 7668.954000000001
$ws.Range("M132").Value = -5138.954000000001
$ws.Range("H134").Value = 5559.017
$ws.Range("I134").Value = 5178.6514
$ws.Range("J134").Value = 6581.25
$ws.Range("K134").Value = 15535.9542
$ws.Range("L134").Value = 19743.75
$ws.Range("M134").Value = -13000.9542
$ws.Range("N134").Value = -24813.75
$ws.Range("H136").Value = 4532.091
$ws.Range("I136").Value = 3726.6365
$ws.Range("K136").Value = 11179.9095
$ws.Range("M136").Value = -8629.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 9299.200000000001
$ws.Range("J31").Value = 9299.200000000001
$ws.Range("L31").Value = 27897.6
$ws.Range("N31").Value = -28473.6
$ws.Range("H137").Value = 4725
$ws.Range("I137").Value = 1880
$ws.Range("J137").Value = 10821.429
$ws.Range("K137").Value = 5640
$ws.Range("L137").Value = 32464.287
$ws.Range("M137").Value = -540
$ws.Range("N137").Value = -42664.287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 127249.5
$ws.Range("I7").Value = 145035.14
$ws.Range("K7").Value = 145035.14
$ws.Range("M7").Value = -144923.14
$ws.Range("H22").Value = 1228.3
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 1228.3
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H68").Value = 2657.0527
$ws.Range("I68").Value = 2623.75
$ws.Range("J68").Value = 2714.1428
$ws.Range("K68").Value = 2623.75
$ws.Range("L68").Value = 2714.1428
$ws.Range("M68").Value = -1874.75
$ws.Range("N68").Value = -4212.1428
$ws.Range("H71").Value = 2657.0527
$ws.Range("I71").Value = 2623.75
$ws.Range("J71").Value = 2714.1428
$ws.Range("K71").Value = 13118.75
$ws.Range("L71").Value = 13570.714
$ws.Range("M71").Value = -9374.75
$ws.Range("N71").Value = -21058.714
$ws.Range("H122").Value = 5950.2188
$ws.Range("I122").Value = 4367.4
$ws.Range("J122").Value = 6669.6816
$ws.Range("K122").Value = 13102.2
$ws.Range("L122").Value = 20009.0448
$ws.Range("M122").Value = -10652.2
$ws.Range("N122").Value = -24909.0448
$ws.Range("H126").Value = 127249.5
$ws.Range("I126").Value = 145035.14
$ws.Range("K126").Value = 435105.42
$ws.Range("M126").Value = -432635.42

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4032.577
$ws.Range("I81").Value = 4243.6665
$ws.Range("J81").Value = 1499.5
$ws.Range("K81").Value = 8487.333000000001
$ws.Range("L81").Value = 2999
$ws.Range("M81").Value = -7426.333000000001
$ws.Range("N81").Value = -5121
$ws.Range("H84").Value = 4032.577
$ws.Range("I84").Value = 4243.6665
$ws.Range("J84").Value = 1499.5
$ws.Range("K84").Value = 42436.665
$ws.Range("L84").Value = 14995
$ws.Range("M84").Value = -37132.665
$ws.Range("N84").Value = -25603

